# Generate Report for Handoff
# Adds a new tracked file (4836d715-d05c-43d8-8054-a2ba861d92df.md) as row 9
# on the Overview sheet and on each per-locale sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$guid = "4836d715-d05c-43d8-8054-a2ba861d92df"
$xlfHash = "59ca9bd87e949d97de0bbc12267b317dedf69698"

$mdName = "$guid.md"
$zhXlfName = "$guid.$xlfHash.zh-cn.xlf"
$deXlfName = "$guid.$xlfHash.de-de.xlf"

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$guid/e2e/$mdName"
$zhHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$xlfHash/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlfName"
$deHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$xlfHash/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlfName"

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1): File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B9").Value = "Ready for handoff"
$wsOverview.Range("C9").Value = "Ready for handoff"
$wsOverview.Range("D9").Value = "2016-03-22 18:43:30"
$wsOverview.Range("D9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A9"), $mdUrl, "", "", $mdName)

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B9").Value = ".md"
$wsZh.Range("C9").Value = "Ready for handoff"
$wsZh.Range("E9").Value = "2016-03-22 18:43:26"
$wsZh.Range("E9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H9").Value = "0001-01-01 00:00:00"
$wsZh.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("J9").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A9"), $mdUrl, "", "", $mdName)
$wsZh.Hyperlinks.Add($wsZh.Range("D9"), $zhHandoffUrl, "", "", $zhXlfName)

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B9").Value = ".md"
$wsDe.Range("C9").Value = "Ready for handoff"
$wsDe.Range("E9").Value = "2016-03-22 18:43:30"
$wsDe.Range("E9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H9").Value = "0001-01-01 00:00:00"
$wsDe.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("J9").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A9"), $mdUrl, "", "", $mdName)
$wsDe.Hyperlinks.Add($wsDe.Range("D9"), $deHandoffUrl, "", "", $deXlfName)
